$d = $word.ActiveDocument

# The "Date" style paragraph ("Version 1.0 vom "2020-10-20"") is being
# replaced by a title image, so remove that whole paragraph (including its
# paragraph mark) from the document.
foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "Date" -or $p.Range.Text -match "Version\s+1\.0\s+vom") {
        $p.Range.Delete()
        break
    }
}
